# Rename the seven "*img" sheets to "img*" and make the last one
# ("eimg" -> "imge") the active/selected tab (previously "holiday" was
# the active tab).

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("himg").Name = "imgh"
$wb.Worksheets.Item("timg").Name = "imgt"
$wb.Worksheets.Item("simg").Name = "imgs"
$wb.Worksheets.Item("gimg").Name = "imgg"
$wb.Worksheets.Item("wimg").Name = "imgw"
$wb.Worksheets.Item("bimg").Name = "imgb"
$wb.Worksheets.Item("eimg").Name = "imge"

$wb.Worksheets.Item("imge").Activate()
